# Remove all the body content between the title heading and the final
# (bookmark-carrying) paragraph: the blank paragraph right after the
# title, the "Produkt" / "Aktiv auktion" / "Bud" / "Budhistorik" sections
# and everything in between, as well as the trailing text of the last
# paragraph -- leaving only the heading paragraph followed by an empty
# paragraph that still carries the _GoBack bookmark.

$d = $word.ActiveDocument

$paragraphs = $d.Paragraphs
$count = $paragraphs.Count

$firstToRemove = $paragraphs.Item(2).Range.Start
$removalEnd = $paragraphs.Item($count - 1).Range.End

if ($removalEnd -gt $firstToRemove) {
    $d.Range($firstToRemove, $removalEnd).Delete()
}

$tail = $paragraphs.Item($paragraphs.Count).Range
$tailTextEnd = $tail.End - 1
if ($tailTextEnd -gt $tail.Start) {
    $d.Range($tail.Start, $tailTextEnd).Delete()
}
